$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 4.913911592104444
$ws.Range("D2").Value = 7.274415023543201
$ws.Range("E2").Value = 12.32293271044325
$ws.Range("F2").Value = 40.13736744714087
$ws.Range("G2").Value = 47.5000846219713
$ws.Range("H2").Value = 18.82453294665168
$ws.Range("J2").Value = 9.970904261452148
$ws.Range("K2").Value = 16.82635729032129
$ws.Range("M2").Value = 18.44911726384874

$ws.Range("C3").Value = 4.753110097562584
$ws.Range("D3").Value = 7.263039525853822
$ws.Range("E3").Value = 12.34152499231859
$ws.Range("F3").Value = 40.18539705760097
$ws.Range("G3").Value = 47.49204344876372
$ws.Range("H3").Value = 18.88320984039449
$ws.Range("J3").Value = 9.998971604615551
$ws.Range("K3").Value = 16.41737856357897
$ws.Range("M3").Value = 18.29791157498233

$ws.Range("C4").Value = 4.653180780207277
$ws.Range("D4").Value = 7.256505783537638
$ws.Range("E4").Value = 12.35460078734784
$ws.Range("F4").Value = 40.22704213976194
$ws.Range("G4").Value = 47.50522067547555
$ws.Range("H4").Value = 18.9237059985647
$ws.Range("J4").Value = 10.0173569626847
$ws.Range("K4").Value = 16.16426269440684
$ws.Range("M4").Value = 18.20781111964726

$ws.Range("C5").Value = 4.612230616187768
$ws.Range("D5").Value = 7.253957932662381
$ws.Range("E5").Value = 12.36034682723063
$ws.Range("F5").Value = 40.24705550817546
$ws.Range("G5").Value = 47.51512519579281
$ws.Range("H5").Value = 18.9413275917865
$ws.Range("J5").Value = 10.02513915199125
$ws.Range("K5").Value = 16.06076647786449
$ws.Range("M5").Value = 18.17181615476086

$ws.Range("C6").Value = 4.605419315565994
$ws.Range("D6").Value = 7.253541829864293
$ws.Range("E6").Value = 12.36132617413957
$ws.Range("F6").Value = 40.25056204822919
$ws.Range("G6").Value = 47.51704303598317
$ws.Range("H6").Value = 18.94432109916081
$ws.Range("J6").Value = 10.02644890475978
$ws.Range("K6").Value = 16.04356452034102
$ws.Range("M6").Value = 18.1658837351116

$ws.Range("C7").Value = 4.652629335955188
$ws.Range("D7").Value = 7.256470956158507
$ws.Range("E7").Value = 12.35467658969623
$ws.Range("F7").Value = 40.22729974697474
$ws.Range("G7").Value = 47.50533592092503
$ws.Range("H7").Value = 18.92393912456517
$ws.Range("J7").Value = 10.0174607413215
$ws.Range("K7").Value = 16.16286811535335
$ws.Range("M7").Value = 18.20732271472034

$ws.Range("C8").Value = 4.858757636694729
$ws.Range("D8").Value = 7.270399982033424
$ws.Range("E8").Value = 12.3289988610102
$ws.Range("F8").Value = 40.15139757145429
$ws.Range("G8").Value = 47.49354146468583
$ws.Range("H8").Value = 18.8438346661911
$ws.Range("J8").Value = 9.980343017401811
$ws.Range("K8").Value = 16.68584704655142
$ws.Range("M8").Value = 18.3964326828233

$ws.Range("C9").Value = 5.250550230160735
$ws.Range("D9").Value = 7.301233402920174
$ws.Range("E9").Value = 12.29181129367658
$ws.Range("F9").Value = 40.09953678904255
$ws.Range("G9").Value = 47.61484030163036
$ws.Range("H9").Value = 18.72240101938159
$ws.Range("J9").Value = 9.916679146049843
$ws.Range("K9").Value = 17.68911965006697
$ws.Range("M9").Value = 18.78753556415611

$ws.Range("C10").Value = 5.527288275435824
$ws.Range("D10").Value = 7.325952600299093
$ws.Range("E10").Value = 12.27250801919332
$ws.Range("F10").Value = 40.1211879737241
$ws.Range("G10").Value = 47.792644045701
$ws.Range("H10").Value = 18.65517950550141
$ws.Range("J10").Value = 9.875445505773232
$ws.Range("K10").Value = 18.40460898814842
$ws.Range("M10").Value = 19.08517665165017

$ws.Range("C11").Value = 5.650133566704874
$ws.Range("D11").Value = 7.337629083751786
$ws.Range("E11").Value = 12.26546508199425
$ws.Range("F11").Value = 40.14410858882562
$ws.Range("G11").Value = 47.89282241053278
$ws.Range("H11").Value = 18.62942796899696
$ws.Range("J11").Value = 9.857885729448894
$ws.Range("K11").Value = 18.72388390553598
$ws.Range("M11").Value = 19.22236689040015

$ws.Range("C12").Value = 5.696168301545296
$ws.Range("D12").Value = 7.342111155372605
$ws.Range("E12").Value = 12.26304775478049
$ws.Range("F12").Value = 40.15467223570721
$ws.Range("G12").Value = 47.93352640070675
$ws.Range("H12").Value = 18.62037486210142
$ws.Range("J12").Value = 9.851408195611
$ws.Range("K12").Value = 18.84377788743193
$ws.Range("M12").Value = 19.27453649941006

$ws.Range("C13").Value = 5.68627611279602
$ws.Range("D13").Value = 7.341143199125483
$ws.Range("E13").Value = 12.26355726973105
$ws.Range("F13").Value = 40.1523133199143
$ws.Range("G13").Value = 47.92463708123568
$ws.Range("H13").Value = 18.62229348241479
$ws.Range("J13").Value = 9.852795603681917
$ws.Range("K13").Value = 18.81800317931677
$ws.Range("M13").Value = 19.26329173844005

$ws.Range("C14").Value = 5.653930821773899
$ws.Range("D14").Value = 7.33799662456832
$ws.Range("E14").Value = 12.26526120480625
$ws.Range("F14").Value = 40.14493988783757
$ws.Range("G14").Value = 47.89611571700995
$ws.Range("H14").Value = 18.62866914326406
$ws.Range("J14").Value = 9.85734937362173
$ws.Range("K14").Value = 18.7337685004498
$ws.Range("M14").Value = 19.22665475735206

$ws.Range("C15").Value = 5.634054062855348
$ws.Range("D15").Value = 7.336077079253148
$ws.Range("E15").Value = 12.26633742102483
$ws.Range("F15").Value = 40.14066891158032
$ws.Range("G15").Value = 47.87900585138793
$ws.Range("H15").Value = 18.63266549869181
$ws.Range("J15").Value = 9.860161078005179
$ws.Range("K15").Value = 18.68203769853159
$ws.Range("M15").Value = 19.20424086293281

$ws.Range("C16").Value = 5.519194709336174
$ws.Range("D16").Value = 7.325198061561511
$ws.Range("E16").Value = 12.27300324685016
$ws.Range("F16").Value = 40.119953514982
$ws.Range("G16").Value = 47.78648502004256
$ws.Range("H16").Value = 18.65695994126361
$ws.Range("J16").Value = 9.876617156532779
$ws.Range("K16").Value = 18.38360861563314
$ws.Range("M16").Value = 19.07624366999429

$ws.Range("C17").Value = 5.44791789544412
$ws.Range("D17").Value = 7.318633535290505
$ws.Range("E17").Value = 12.27753753555056
$ws.Range("F17").Value = 40.11059700359855
$ws.Range("G17").Value = 47.73466568985421
$ws.Range("H17").Value = 18.67310337436665
$ws.Range("J17").Value = 9.887018989644391
$ws.Range("K17").Value = 18.19885718595976
$ws.Range("M17").Value = 18.99815274637699

$ws.Range("C18").Value = 5.406637254967663
$ws.Range("D18").Value = 7.314898516552002
$ws.Range("E18").Value = 12.28030916737591
$ws.Range("F18").Value = 40.10644558843149
$ws.Range("G18").Value = 47.70667735494031
$ws.Range("H18").Value = 18.68284283232789
$ws.Range("J18").Value = 9.893114588012804
$ws.Range("K18").Value = 18.09201539018794
$ws.Range("M18").Value = 18.95340791530388

$ws.Range("C19").Value = 5.392613041322395
$ws.Range("D19").Value = 7.313640946505758
$ws.Range("E19").Value = 12.28127570523068
$ws.Range("F19").Value = 40.10525110968101
$ws.Range("G19").Value = 47.69751311490707
$ws.Range("H19").Value = 18.68621832060062
$ws.Range("J19").Value = 9.895197823456394
$ws.Range("K19").Value = 18.05574507463373
$ws.Range("M19").Value = 18.93828864102075

$ws.Range("C20").Value = 5.455535189454966
$ws.Range("D20").Value = 7.319328137251882
$ws.Range("E20").Value = 12.27703791965266
$ws.Range("F20").Value = 40.11146566049008
$ws.Range("G20").Value = 47.73999394407016
$ws.Range("H20").Value = 18.67133784138975
$ws.Range("J20").Value = 9.885900030503175
$ws.Range("K20").Value = 18.21858494828372
$ws.Range("M20").Value = 19.00644824274029

$ws.Range("C21").Value = 5.663444900638521
$ws.Range("D21").Value = 7.33891922254037
$ws.Range("E21").Value = 12.26475394394398
$ws.Range("F21").Value = 40.14705448404493
$ws.Range("G21").Value = 47.90441806585511
$ws.Range("H21").Value = 18.62677746909484
$ws.Range("J21").Value = 9.856007155619126
$ws.Range("K21").Value = 18.75853851855523
$ws.Range("M21").Value = 19.23741030419906

$ws.Range("C22").Value = 5.79648617315174
$ws.Range("D22").Value = 7.352074648579931
$ws.Range("E22").Value = 12.25818082157768
$ws.Range("F22").Value = 40.18129647553235
$ws.Range("G22").Value = 48.02801004179174
$ws.Range("H22").Value = 18.60172755223188
$ws.Range("J22").Value = 9.837472713587669
$ws.Range("K22").Value = 19.1055036612883
$ws.Range("M22").Value = 19.38961367209898

$ws.Range("C23").Value = 5.725753391857602
$ws.Range("D23").Value = 7.345021730542414
$ws.Range("E23").Value = 12.2615559759813
$ws.Range("F23").Value = 40.16201502493491
$ws.Range("G23").Value = 47.96057392503233
$ws.Range("H23").Value = 18.61472315875545
$ws.Range("J23").Value = 9.847273263531887
$ws.Range("K23").Value = 18.92090006124032
$ws.Range("M23").Value = 19.30827781777246

$ws.Range("C24").Value = 5.452092349981953
$ws.Range("D24").Value = 7.319013986070948
$ws.Range("E24").Value = 12.27726328262959
$ws.Range("F24").Value = 40.11106911666239
$ws.Range("G24").Value = 47.73757942287208
$ws.Range("H24").Value = 18.67213461131392
$ws.Range("J24").Value = 9.88640555230938
$ws.Range("K24").Value = 18.20966796703437
$ws.Range("M24").Value = 19.00269737650932

$ws.Range("C25").Value = 5.146276889245844
$ws.Range("D25").Value = 7.292523475461615
$ws.Range("E25").Value = 12.30046263145614
$ws.Range("F25").Value = 40.10311145436739
$ws.Range("G25").Value = 47.56647854406846
$ws.Range("H25").Value = 18.75140770568251
$ws.Range("J25").Value = 9.932927436450914
$ws.Range("K25").Value = 17.42095556278504
$ws.Range("M25").Value = 18.67977190516985

Write-Output "Updated loading_percent values for rows 2-25"